$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 with the new "Sauteruz" parish data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5070000000
$ws.Range("C6").Value = 5070
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = "Sauteruz"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "P"

# Update the selection to E9
$ws.Range("E9").Select()
